$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-26 Friday", "2024-04-27 Saturday"),
    @("77×24=1848", "57×44=2508"),
    @("58×14=812", "38×24=912"),
    @("64×43=2752", "66×36=2376"),
    @("90×76=6840", "33×70=2310"),
    @("78×42=3276", "36×35=1260"),
    @("13×66=858", "40×26=1040"),
    @("51×58=2958", "95×85=8075"),
    @("77×82=6314", "83×99=8217"),
    @("47×66=3102", "54×68=3672"),
    @("18×15=270", "91×67=6097"),
    @("97×51=4947", "30×64=1920"),
    @("70×45=3150", "96×81=7776"),
    @("97×53=5141", "44×11=484"),
    @("42×85=3570", "89×57=5073"),
    @("51×18=918", "20×66=1320"),
    @("65×78=5070", "98×74=7252"),
    @("96×73=7008", "62×51=3162"),
    @("17×85=1445", "13×17=221"),
    @("51×24=1224", "77×39=3003"),
    @("28×53=1484", "70×49=3430"),
    @("12×95=1140", "50×87=4350"),
    @("81×68=5508", "57×35=1995"),
    @("48×48=2304", "19×17=323"),
    @("24×63=1512", "37×30=1110"),
    @("69×21=1449", "85×37=3145")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
